# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the daily price rows (2-11) for
# "Mora" at Vega Central Mapocho de Santiago: each row's Fecha/Volumen/
# Precio minimo/Precio maximo/Precio promedio ponderado/Origen/Precio $/Kg
# values get redistributed among the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move (D, M, N, O, P, R, S)
# for data rows 2..11, before writing anything, since values are reused
# across rows (a pure permutation of existing row data).
$snapshot = @{}
for ($r = 2; $r -le 11; $r++) {
    $row = @{
        D = $ws.Cells.Item($r, 4).Value()
        M = $ws.Cells.Item($r, 13).Value()
        N = $ws.Cells.Item($r, 14).Value()
        O = $ws.Cells.Item($r, 15).Value()
        P = $ws.Cells.Item($r, 16).Value()
        R = $ws.Cells.Item($r, 18).Value()
        S = $ws.Cells.Item($r, 19).Value()
    }
    $snapshot[$r] = $row
}

# Mapping of destination row -> source row (which row's data it now holds)
$mapping = @{
    2  = 9
    3  = 10
    4  = 2
    5  = 8
    6  = 7
    7  = 3
    8  = 5
    9  = 4
    10 = 11
    11 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $data.D
    $ws.Cells.Item($destRow, 13).Value = $data.M
    $ws.Cells.Item($destRow, 14).Value = $data.N
    $ws.Cells.Item($destRow, 15).Value = $data.O
    $ws.Cells.Item($destRow, 16).Value = $data.P
    $ws.Cells.Item($destRow, 18).Value = $data.R
    $ws.Cells.Item($destRow, 19).Value = $data.S
}
